$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A1").Select()
Write-Host "ok"
